$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "theta_threshold_range" row (row 5). This shifts the
# "pie_threshold_range" row up from row 6 to row 5, and Excel/the engine
# drops the now-unused "theta_threshold_range" shared string automatically.
$ws.Rows("5").Delete()

# Update the Min values for the remaining parameter rows.
$ws.Range("B2").Value = 5.4   # alpha_distance_range Min: 3.8 -> 5.4
$ws.Range("B3").Value = 5.7   # beta_distance_range Min: 5.3 -> 5.7
$ws.Range("B4").Value = 0.7   # ratio_threshold_range Min: 0.8 -> 0.7

# Row 5 now holds the former pie_threshold_range row (previously row 6).
$ws.Range("B5").Value = 0     # pie_threshold_range Min stays 0
$ws.Range("C5").Value = 15    # pie_threshold_range Max: 20 -> 15

# Match the author's final selection (active cell C4).
$ws.Range("C4").Select()
